# Edit the "How Does a Bike Share Navigate Speedy Success" presentation:
# update the closing SmartArt slide (last slide) text content and
# reposition/resize the SmartArt graphic frame.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item($p.Slides.Count)

# Locate the SmartArt graphic frame shape on the slide.
$smartArtShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasSmartArt) {
        $smartArtShape = $candidate
    }
}

$sa = $smartArtShape.SmartArt

# Update the three text nodes (order matches the underlying diagram data).
$sa.AllNodes.Item(1).TextFrame2.TextRange.Text = "Thanks for listening"
$sa.AllNodes.Item(2).TextFrame2.TextRange.Text = "Welcome to ask questions"
$sa.AllNodes.Item(3).TextFrame2.TextRange.Text = "If the answer to this question is incomplete, please let me know."

# Reposition / resize the SmartArt graphic frame itself (EMU values chosen
# so the point -> EMU round trip lands exactly on the target offsets).
$smartArtShape.Left = 0.0
$smartArtShape.Top = 124.09095001220703
$smartArtShape.Width = 960.0
$smartArtShape.Height = 415.9090576171875
